$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab to match the new file name
$ws.Name = "2023-pmhc-outcome-measure-stand"

# Update the report title date (A1) to reflect the 30 June reporting date
$ws.Range("A1").Value = "Report Z1 Measure standard deviations for National; 30/06/2023"

# Update the "Generated on" timestamp (A10)
$ws.Range("A10").Value = "Generated on 28/09/2023 12:13:02 AEST"

# Update the recalculated standard deviation values
$ws.Range("C4").Value = 6.1638584959999996
$ws.Range("C5").Value = 5.9116280210000003
$ws.Range("C6").Value = 5.7441892650000002
$ws.Range("C7").Value = 4.943526233
$ws.Range("C8").Value = 8.7272310110000006

# Restore default Excel page margins
# (PageSetup margins are expressed in points; OOXML pageMargins are in
# inches, so multiply the target inch values by 72 points/inch.)
$ws.PageSetup.LeftMargin = 0.75 * 72
$ws.PageSetup.RightMargin = 0.75 * 72
$ws.PageSetup.TopMargin = 1 * 72
$ws.PageSetup.BottomMargin = 1 * 72
$ws.PageSetup.HeaderMargin = 0.5 * 72
$ws.PageSetup.FooterMargin = 0.5 * 72
